$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# "Status" column text changed from "Ready for handoff" to "Handed back: in sync with en-US"
# (shared by Overview!E2/F2 and the "Status" column on both language sheets)
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"

# zh-cn: Latest Handback DateTime refreshed, Error Detail cleared (handback now in sync)
$zhcn.Range("K2").Value = "2016-09-06 09:07:28"
$zhcn.Range("P2").Value = ""

# de-de: Latest Handback DateTime refreshed, Error Detail cleared (handback now in sync)
$dede.Range("K2").Value = "2016-09-06 09:07:36"
$dede.Range("P2").Value = ""

# Column width adjustments (report regenerated with wider Status/Error Detail columns)
$overview.Columns.Item(5).ColumnWidth = 29.09
$overview.Columns.Item(6).ColumnWidth = 29.09

$zhcn.Columns.Item(3).ColumnWidth = 29.09
$zhcn.Columns.Item(16).ColumnWidth = 12.91

$dede.Columns.Item(3).ColumnWidth = 29.09
$dede.Columns.Item(16).ColumnWidth = 12.91
